$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlUp = -4162 : find the last used row in column A, then target the row right after it
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = 10002
$ws.Cells.Item($newRow, 2).Value = 110032
$ws.Cells.Item($newRow, 3).Value = "eng"
$ws.Cells.Item($newRow, 4).Value = $true
$ws.Cells.Item($newRow, 5).Value = "superadmin"
$ws.Cells.Item($newRow, 6).Value = "now()"

# Match the author's final selection / scroll position
$ws.Range("C31").Select()
$excel.ActiveWindow.ScrollRow = 22
